# "add date of retrieval and fix title-language error"
#
# The underlying data fix: rows 27, 29, 32 and 37 were tagged as a "T"
# (German-titled) row in column A but should have been tagged "F" (i.e. not
# a German-specific title) like their sibling rows - this is the
# "title-language error" referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "F"
$ws.Range("A29").Value = "F"
$ws.Range("A32").Value = "F"
$ws.Range("A37").Value = "F"

# The author also scrolled/re-selected the sheet while reviewing the fix
# (new cursor position at E18, with row 10 pinned to the top of the view).
$ws.Range("E18").Select()
